# Update gh-pages output (generated at 456a3b4)
$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (sheet1) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = "不可售"
$ws1.Range("F4").Value = 295
$ws1.Range("F5").Value = 858
$ws1.Range("F6").Value = 18
$ws1.Range("F7").Value = 306
$ws1.Range("F8").Value = 8776
$ws1.Range("F9").Value = 77
$ws1.Range("F11").Value = 129
$ws1.Range("F12").Value = 115
$ws1.Range("F13").Value = 8
$ws1.Range("F17").Value = 28
$ws1.Range("F18").Value = 259
$ws1.Range("F19").Value = 744
$ws1.Range("F20").Value = 36

# ---- Sheet "全部类型" (sheet4) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 269
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F4").Value = 295
$ws4.Range("F5").Value = 858
$ws4.Range("F6").Value = 18
$ws4.Range("F7").Value = 306
$ws4.Range("F8").Value = 8776
$ws4.Range("F9").Value = 77
$ws4.Range("F11").Value = 129
$ws4.Range("F12").Value = 115
$ws4.Range("F13").Value = 8
$ws4.Range("F17").Value = 28
$ws4.Range("F18").Value = 259
$ws4.Range("F19").Value = 744
$ws4.Range("F20").Value = 36
